# Apply cryptos list update (Sun Nov 24 09:43:09 UTC 2024)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '97.942.65'
$ws.Range('E2').Value = '  -0.55%  '
$ws.Range('D3').Value = '3.395.92'
$ws.Range('E3').Value = '  +0.85%  '
$ws.Range('E4').Value = '  +0.04%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '253.50'
$ws.Range('E5').Value = '  -1.80%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '676.60'
$ws.Range('E6').Value = '  +1.40%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '1.43'
$ws.Range('E7').Value = '  -6.69%  '
$ws.Range('E8').Value = '  -7.78%  '
$ws.Range('E9').Value = '  -3.30%  '
$ws.Range('E10').Value = '  +0.00%  '
$ws.Range('D11').Value = '3.395.14'
$ws.Range('E11').Value = '  +0.91%  '
$ws.Range('E12').Value = '  +1.08%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '41.46'
$ws.Range('E13').Value = '  -1.42%  '
$ws.Range('E14').Value = '  +10.34%  '
$ws.Range('D15').Value = '97.636.05'
$ws.Range('E15').Value = '  -0.67%  '
$ws.Range('E16').Value = '  -2.72%  '
$ws.Range('D17').Value = '4.023.97'
$ws.Range('E17').Value = '  +0.64%  '
$ws.Range('E18').Value = '  +16.22%  '
$ws.Range('D19').Value = '3.381.97'
$ws.Range('E19').Value = '  +0.41%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.565'
$ws.Range('E20').Value = '  +29.20%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '17.30'
$ws.Range('E21').Value = '  +2.76%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '10.97'
$ws.Range('E22').Value = '  +3.98%  '
$ws.Range('E23').Value = '  -4.38%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '505.32'
$ws.Range('E24').Value = '  -4.75%  '
$ws.Range('E25').Value = '  -6.82%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '6.51'
$ws.Range('E26').Value = '  +4.09%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '99.01'
$ws.Range('E27').Value = '  -3.63%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '12.52'
$ws.Range('E28').Value = '  -0.71%  '
$ws.Range('D29').Value = '3.580.31'
$ws.Range('E29').Value = '  +0.98%  '
$ws.Range('E30').Value = '  -0.31%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '11.45'
$ws.Range('E31').Value = '  +3.65%  '
$ws.Range('E32').Value = '  +0.21%  '
$ws.Range('E33').Value = '  +1.64%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '2.58'
$ws.Range('E34').Value = '  +21.02%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.00'
$ws.Range('E35').Value = '  -0.19%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.566'
$ws.Range('E36').Value = '  +3.53%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '29.32'
$ws.Range('E37').Value = '  -0.80%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '1.50'
$ws.Range('E38').Value = '  +10.75%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '7.86'
$ws.Range('E39').Value = '  -0.38%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '529.47'
$ws.Range('E40').Value = '  -0.21%  '
$ws.Range('E41').Value = '  -4.24%  '
$ws.Range('E42').Value = '  +0.00%  '
$ws.Range('E43').Value = '  +0.05%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.865'
$ws.Range('E44').Value = '  +3.52%  '
$ws.Range('B45').Value = 'MantraDAO'
$ws.Range('C45').Value = 'https://coinranking.com/coin/cTdD8lD-6+mantradao-om'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '3.76'
$ws.Range('E45').Value = '  +0.00%  '
$ws.Range('B46').Value = 'VeChain'
$ws.Range('C46').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.0427'
$ws.Range('E46').Value = '  -2.08%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '8.86'
$ws.Range('E47').Value = '  +11.49%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.72'
$ws.Range('E48').Value = '  +12.02%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '5.71'
$ws.Range('E49').Value = '  +11.54%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '55.49'
$ws.Range('E50').Value = '  +10.13%  '
$ws.Range('E51').Value = '  -7.22%  '
